$wb = $excel.ActiveWorkbook
$wsYearly = $wb.Worksheets.Item("Yearly")
$wsAllTime = $wb.Worksheets.Item("All Time")

# --- Data edit: Yearly!M14 changed (other touched cells are dependent
# formulas that recalculate automatically: O14, M15, O15 on Yearly and
# G8, I8, G46, I46 on All Time, which reference Yearly totals). ---
$wsYearly.Range("M14").Value = 105.36

# --- View/selection state ---
# Yearly sheet: selection moves from N15 to N14. Select it first so the
# final active/selected sheet ends up being "All Time" (matching the
# unchanged tabSelected/activeTab state in the workbook).
$wsYearly.Range("N14").Select()

# All Time sheet: becomes (stays) the active sheet, scrolled so row 19 is
# at the top, with the selection moved from M29 to L19.
$wsAllTime.Activate()
$wsAllTime.Range("L19").Select()
$excel.ActiveWindow.ScrollRow = 19
$excel.ActiveWindow.ScrollColumn = 1
